# Update "Gal-Gpr151" LR-pairs sheet with refreshed TPM-derived NATMI output.
# The recomputation adds "ECs" as a sending cluster (previously it only
# appeared as a target cluster) and refreshes all derived-specificity /
# expression-weight figures for every Sending-cluster x Target-cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Pipe-delimited row data: for each of the 20 columns, a type flag (S=string,
# N=number) followed by the literal value. 16 data rows (rows 2-17).
$rowData = @"
S|ECs|S|Gal|S|Gpr151|S|ECs|N|2|N|0.6666666666666666|N|1.061714666666667|N|3.185144|N|0.09827129933343294|N|0.09827129933343293|N|3|N|1|N|0.336217|N|1.008651|N|0.1376522887391218|N|0.1376522887391219|N|0.3569665200826667|N|3.212698680744|N|0.01352726927061438|N|0.01352726927061438
S|ECs|S|Gal|S|Gpr151|S|FAPs|N|2|N|0.6666666666666666|N|1.061714666666667|N|3.185144|N|0.09827129933343294|N|0.09827129933343293|N|3|N|1|N|0.5236423333333333|N|1.570927|N|0.2143870347544219|N|0.214387034754422|N|0.5559587453875556|N|5.003628708488|N|0.02106809246555889|N|0.02106809246555889
S|ECs|S|Gal|S|Gpr151|S|MuSCs|N|2|N|0.6666666666666666|N|1.061714666666667|N|3.185144|N|0.09827129933343294|N|0.09827129933343293|N|3|N|1|N|1.197380666666667|N|3.592142|N|0.4902256258863835|N|0.4902256258863836|N|1.271276615383111|N|11.441489538448|N|0.0481751092224003|N|0.04817510922240031
S|ECs|S|Gal|S|Gpr151|S|Resolving-Mac|N|2|N|0.6666666666666666|N|1.061714666666667|N|3.185144|N|0.09827129933343294|N|0.09827129933343293|N|3|N|1|N|0.3852693333333333|N|1.155808|N|0.1577350506200726|N|0.1577350506200727|N|0.4090461018168889|N|3.681414916352|N|0.01550082837485936|N|0.01550082837485936
S|FAPs|S|Gal|S|Gpr151|S|ECs|N|3|N|1|N|1.642758|N|4.928274|N|0.1520521174085614|N|0.1520521174085614|N|3|N|1|N|0.336217|N|1.008651|N|0.1376522887391218|N|0.1376522887391219|N|0.552323166486|N|4.970908498374|N|0.02093032196891814|N|0.02093032196891815
S|FAPs|S|Gal|S|Gpr151|S|FAPs|N|3|N|1|N|1.642758|N|4.928274|N|0.1520521174085614|N|0.1520521174085614|N|3|N|1|N|0.5236423333333333|N|1.570927|N|0.2143870347544219|N|0.214387034754422|N|0.860217632222|N|7.741958689998|N|0.03259800257935268|N|0.0325980025793527
S|FAPs|S|Gal|S|Gpr151|S|MuSCs|N|3|N|1|N|1.642758|N|4.928274|N|0.1520521174085614|N|0.1520521174085614|N|3|N|1|N|1.197380666666667|N|3.592142|N|0.4902256258863835|N|0.4902256258863836|N|1.967006669212|N|17.703060022908|N|0.07453984442396186|N|0.07453984442396187
S|FAPs|S|Gal|S|Gpr151|S|Resolving-Mac|N|3|N|1|N|1.642758|N|4.928274|N|0.1520521174085614|N|0.1520521174085614|N|3|N|1|N|0.3852693333333333|N|1.155808|N|0.1577350506200726|N|0.1577350506200727|N|0.6329042794879999|N|5.696138515392|N|0.02398394843632866|N|0.02398394843632867
S|MuSCs|S|Gal|S|Gpr151|S|ECs|N|3|N|1|N|7.792831666666667|N|23.378495|N|0.7212970842480482|N|0.7212970842480481|N|3|N|1|N|0.336217|N|1.008651|N|0.1376522887391218|N|0.1376522887391219|N|2.620082484471667|N|23.580742360245|N|0.09928819450759901|N|0.09928819450759901
S|MuSCs|S|Gal|S|Gpr151|S|FAPs|N|3|N|1|N|7.792831666666667|N|23.378495|N|0.7212970842480482|N|0.7212970842480481|N|3|N|1|N|0.5236423333333333|N|1.570927|N|0.2143870347544219|N|0.214387034754422|N|4.080656557207223|N|36.725909014865|N|0.1546367430689495|N|0.1546367430689495
S|MuSCs|S|Gal|S|Gpr151|S|MuSCs|N|3|N|1|N|7.792831666666667|N|23.378495|N|0.7212970842480482|N|0.7212970842480481|N|3|N|1|N|1.197380666666667|N|3.592142|N|0.4902256258863835|N|0.4902256258863836|N|9.330985976254444|N|83.97887378628999|N|0.3535983145755229|N|0.3535983145755229
S|MuSCs|S|Gal|S|Gpr151|S|Resolving-Mac|N|3|N|1|N|7.792831666666667|N|23.378495|N|0.7212970842480482|N|0.7212970842480481|N|3|N|1|N|0.3852693333333333|N|1.155808|N|0.1577350506200726|N|0.1577350506200727|N|3.002339060995555|N|27.02105154896|N|0.1137738320959767|N|0.1137738320959767
S|Resolving-Mac|S|Gal|S|Gpr151|S|ECs|N|3|N|1|N|0.3066096666666667|N|0.919829|N|0.02837949900995756|N|0.02837949900995756|N|3|N|1|N|0.336217|N|1.008651|N|0.1376522887391218|N|0.1376522887391219|N|0.1030873822976667|N|0.927786440679|N|0.003906502991990301|N|0.003906502991990302
S|Resolving-Mac|S|Gal|S|Gpr151|S|FAPs|N|3|N|1|N|0.3066096666666667|N|0.919829|N|0.02837949900995756|N|0.02837949900995756|N|3|N|1|N|0.5236423333333333|N|1.570927|N|0.2143870347544219|N|0.214387034754422|N|0.1605538012758889|N|1.444984211483|N|0.006084196640560854|N|0.006084196640560856
S|Resolving-Mac|S|Gal|S|Gpr151|S|MuSCs|N|3|N|1|N|0.3066096666666667|N|0.919829|N|0.02837949900995756|N|0.02837949900995756|N|3|N|1|N|1.197380666666667|N|3.592142|N|0.4902256258863835|N|0.4902256258863836|N|0.3671284870797778|N|3.304156383718|N|0.01391235766449845|N|0.01391235766449845
S|Resolving-Mac|S|Gal|S|Gpr151|S|Resolving-Mac|N|3|N|1|N|0.3066096666666667|N|0.919829|N|0.02837949900995756|N|0.02837949900995756|N|3|N|1|N|0.3852693333333333|N|1.155808|N|0.1577350506200726|N|0.1577350506200727|N|0.1181273018702222|N|1.063145716832|N|0.004476441712907958|N|0.00447644171290796
"@

$lines = $rowData -split "`r?`n"
$rowNum = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $type = $parts[$i * 2]
        $val = $parts[$i * 2 + 1]
        $cellRef = $cols[$i] + $rowNum
        if ($type -eq "S") {
            $ws.Range($cellRef).Value = $val
        } else {
            $ws.Range($cellRef).Value = [double]$val
        }
    }
    $rowNum += 1
}

Write-Host "Updated rows 2..$($rowNum - 1) on sheet '$($ws.Name)'."
Write-Host "Used range now:" $ws.UsedRange.Address()
